$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is plain numeric-looking text (e.g. "571.64")
# must be forced to Text format first, otherwise Excel auto-converts the
# assignment into a real number and the string representation (and the
# original "t=inlineStr"/shared-string semantics) is lost. We flip the
# format to Text, write the value, then reset the style back to "Normal"
# so no stray explicit style index is left behind on the cell.

$ws.Range("D2").Value = "63.359.26"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "2.547.69"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "2.541.02"
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.84%  "
$ws.Range("D15").Value = "3.003.53"
$ws.Range("E15").Value = "  +4.79%  "
$ws.Range("D16").Value = "63.274.08"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "2.549.47"
$ws.Range("E18").Value = "  +4.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.50%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.80%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +11.35%  "
$ws.Range("E30").Value = "  +12.12%  "
$ws.Range("D31").Value = "0.0₃0838"
$ws.Range("E31").Value = "  +6.85%  "
$ws.Range("E32").Value = "  +3.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("E34").Value = "  +8.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "420.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.408"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "156.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.16%  "
$ws.Range("E44").Value = "  +4.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.609"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0967"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.59%  "
$ws.Range("E51").Value = "  +7.37%  "
